$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: 2021年 ---
# Copy formatting from an existing year-label cell (A2) so the new label
# reuses the same bold/centered/bordered style as the rest of column A.
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A7").Value = "2021年"

$ws.Range("B7").Value = 102.1
$ws.Range("C7").Value = 101.1
$ws.Range("D7").Value = 102.7
$ws.Range("E7").Value = 100.5

# F7 has no data point (blank text cell, matches the empty cells already
# present elsewhere in the sheet, e.g. F6). Entering then clearing the
# formatting keeps the cell present with an empty text value but no
# special quote-prefix formatting.
$ws.Range("F7").Value = "'"
$ws.Range("F7").ClearFormats()

# --- Row 8: 2022年 (only the C column has a value so far) ---
$ws.Range("A2").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A8").Value = "2022年"

$ws.Range("B8").Value = "'"
$ws.Range("B8").ClearFormats()

$ws.Range("C8").Value = 101.3

$ws.Range("D8").Value = "'"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").Value = "'"
$ws.Range("E8").ClearFormats()

$ws.Range("F8").Value = "'"
$ws.Range("F8").ClearFormats()
